# Prepare publication 3.1.0 (#114)
# - Bumps the "Date" metadata value
# - Adds a new "Note that FHIR strings SHALL NOT exceed 1MB in size" comment
#   for the Coding.code element rows (14 and 15) on the Elements sheet
# - Sets the Condition(s) column (AI) to "ele-1" for several element rows
# - Normalizes a couple of "n/a" Mapping values to "N/A" on the Elements sheet

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the publication Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-10-02T08:34:04+00:00"

# --- Elements sheet: apply the row level updates ---
$elements = $wb.Worksheets.Item("Elements")

# Condition(s) column (AI) -> "ele-1" for rows 4,6,8,9,11,12,13,14,15,16,17
$conditionRows = @(4, 6, 8, 9, 11, 12, 13, 14, 15, 16, 17)
foreach ($r in $conditionRows) {
    $elements.Range("AI$r").Value = "ele-1`n"
}

# Mapping: RIM Mapping column (AK) -> "N/A" (was "n/a") for rows 8 and 11
$elements.Range("AK8").Value = "N/A"
$elements.Range("AK11").Value = "N/A"

# Comments column (N) -> new note for the Coding.code / Coding.display rows (14, 15)
$elements.Range("N14").Value = "Note that FHIR strings SHALL NOT exceed 1MB in size"
$elements.Range("N15").Value = "Note that FHIR strings SHALL NOT exceed 1MB in size"

$wb.Save()
